$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Omar Moore'
$ws.Cells.Item(2, 4).Value = '255 Sheila Forks' + [char]10 + 'Port Lisa, MT 30599'
$ws.Cells.Item(3, 1).Value = 'Courtney Lopez'
$ws.Cells.Item(3, 4).Value = '88657 Bridges Ridges Suite 908' + [char]10 + 'New James, FM 12406'
$ws.Cells.Item(4, 1).Value = 'April Acevedo'
$ws.Cells.Item(4, 4).Value = '5782 Karen Freeway Apt. 121' + [char]10 + 'East Michael, MN 93036'
$ws.Cells.Item(5, 1).Value = 'Larry Fitzgerald'
$ws.Cells.Item(5, 4).Value = '05018 Jonathan Causeway Suite 007' + [char]10 + 'Robertberg, PA 76718'
$ws.Cells.Item(6, 1).Value = 'Julia Turner'
$ws.Cells.Item(6, 4).Value = '557 Deborah Stravenue Suite 512' + [char]10 + 'Hunterland, PR 10743'
$ws.Cells.Item(7, 1).Value = 'Eric Gibbs'
$ws.Cells.Item(7, 4).Value = '42438 Choi Hill' + [char]10 + 'Anthonyfurt, AL 77481'
$ws.Cells.Item(8, 1).Value = 'Shannon Chen MD'
$ws.Cells.Item(8, 4).Value = '7581 Ricky Ville Suite 054' + [char]10 + 'Edwardstown, CA 70682'
$ws.Cells.Item(9, 1).Value = 'Jeffrey Dean'
$ws.Cells.Item(9, 4).Value = '50067 Jennifer Keys Apt. 550' + [char]10 + 'Ryanfort, MA 75532'
$ws.Cells.Item(10, 1).Value = 'Rebecca Mcdaniel'
$ws.Cells.Item(10, 4).Value = '313 Ronald Throughway Suite 824' + [char]10 + 'Dennisside, OR 95517'
$ws.Cells.Item(11, 1).Value = 'Keith Browning'
$ws.Cells.Item(11, 4).Value = '4141 Reed Village Apt. 499' + [char]10 + 'Brianbury, SD 77602'
$ws.Cells.Item(12, 1).Value = 'Laura Sanders'
$ws.Cells.Item(12, 4).Value = 'PSC 1945, Box 0739' + [char]10 + 'APO AA 83783'
$ws.Cells.Item(13, 1).Value = 'Rachel Perry'
$ws.Cells.Item(13, 4).Value = '19535 Joel Village' + [char]10 + 'Brownside, WV 45373'
$ws.Cells.Item(14, 1).Value = 'Thomas Ray'
$ws.Cells.Item(14, 4).Value = '154 Sonia Loaf' + [char]10 + 'Lake Danieltown, IA 33784'
$ws.Cells.Item(15, 1).Value = 'Lisa Jensen'
$ws.Cells.Item(15, 4).Value = '69188 Thomas Centers' + [char]10 + 'East Kimberlychester, VI 29574'
$ws.Cells.Item(16, 1).Value = 'Shelby Zamora'
$ws.Cells.Item(16, 4).Value = '05547 Luna Crossroad' + [char]10 + 'Stonefurt, IN 30064'
$ws.Cells.Item(17, 1).Value = 'Amanda Tate'
$ws.Cells.Item(17, 4).Value = '746 Brandon Locks Apt. 876' + [char]10 + 'Gabrielaville, OK 12010'
$ws.Cells.Item(18, 1).Value = 'Lee Faulkner'
$ws.Cells.Item(18, 4).Value = '242 Martin Radial' + [char]10 + 'Greenland, DC 64953'
$ws.Cells.Item(19, 1).Value = 'John Patterson'
$ws.Cells.Item(19, 4).Value = '5134 David Via' + [char]10 + 'Lake Evanbury, NV 59386'
$ws.Cells.Item(20, 1).Value = 'Deanna Allen'
$ws.Cells.Item(20, 4).Value = '64830 Reyes Isle' + [char]10 + 'North Chad, PA 72164'
$ws.Cells.Item(21, 1).Value = 'Wanda Ochoa'
$ws.Cells.Item(21, 4).Value = '53668 Alisha Village' + [char]10 + 'Christopherton, MI 28909'
$ws.Cells.Item(22, 1).Value = 'Hayden Patterson'
$ws.Cells.Item(22, 4).Value = '55868 Mary Falls Suite 905' + [char]10 + 'Kimberlyland, TX 41676'
$ws.Cells.Item(23, 1).Value = 'Chris Weber'
$ws.Cells.Item(23, 4).Value = '6349 Marissa Pass Suite 370' + [char]10 + 'Sethton, WY 82983'
$ws.Cells.Item(24, 1).Value = 'Todd Martinez'
$ws.Cells.Item(24, 4).Value = 'PSC 3488, Box 2229' + [char]10 + 'APO AA 43527'
$ws.Cells.Item(25, 1).Value = 'Rodney Rangel'
$ws.Cells.Item(25, 4).Value = '28088 John Mall' + [char]10 + 'Terrishire, IA 76401'
$ws.Cells.Item(26, 1).Value = 'Ashley Brown'
$ws.Cells.Item(26, 4).Value = '5591 Baker Groves Suite 620' + [char]10 + 'Lorimouth, AL 54896'
$ws.Cells.Item(27, 1).Value = 'Dr. Mark Delgado'
$ws.Cells.Item(27, 4).Value = '5612 Jeff Roads Suite 082' + [char]10 + 'Khanville, RI 13665'
$ws.Cells.Item(28, 1).Value = 'Michael Williams'
$ws.Cells.Item(28, 4).Value = '55013 Catherine Meadow Apt. 040' + [char]10 + 'North Maurice, SD 76429'
$ws.Cells.Item(29, 1).Value = 'Wendy Malone'
$ws.Cells.Item(29, 4).Value = '13621 Michael Harbor' + [char]10 + 'South Bettyshire, NH 58398'
$ws.Cells.Item(30, 1).Value = 'Brittany Garcia'
$ws.Cells.Item(30, 4).Value = '9982 Lawrence Manors Apt. 459' + [char]10 + 'Davidstad, MN 21595'
$ws.Cells.Item(31, 1).Value = 'Brian Lucas'
$ws.Cells.Item(31, 4).Value = '531 Reilly Trail' + [char]10 + 'Rayville, TN 94845'
